# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (fund-holding detail, like the existing
# "2021-Q4" sheet) positioned between "2021-Q4" and "总计", and updates the
# "总计" (totals) sheet with a new summary row for "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item(1)          # "2021-Q4"

$ws = $wb.Worksheets.Add($null, $wsQ4)
$ws.Name = "2022-Q1"

# Look the totals sheet up by name (NOT a reference captured before the
# insertion above) - Worksheets.Item(N) is a live positional lookup, so a
# handle grabbed before the new sheet shifted everything down would now
# silently point at "2022-Q1" instead of "总计".
$wsTotal = $wb.Worksheets.Item("总计")

# --- header row (bold / centered, same look as the "2021-Q4" header) ------
$ws.Range("B1:H1").NumberFormat = "@"
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$wsQ4.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# --- data rows --------------------------------------------------------------
$ws.Range("B2:G3").NumberFormat = "@"

$ws.Range("B2").Value = "000314"
$ws.Range("C2").Value = "招商瑞丰灵活配置混合A"
$ws.Range("D2").Value = "14.01"
$ws.Range("E2").Value = "40.96"
$ws.Range("F2").Value = "1.41"
$ws.Range("G2").Value = "0.1975"
$ws.Range("H2").Value = 10

$ws.Range("B3").Value = "002017"
$ws.Range("C3").Value = "招商瑞丰灵活配置混合C"
$ws.Range("D3").Value = "4.66"
$ws.Range("E3").Value = "40.96"
$ws.Range("F3").Value = "1.41"
$ws.Range("G3").Value = "0.0657"
$ws.Range("H3").Value = 10

$ws.Range("B2:G3").Style = "Normal"

# --- leading index column (A2:A3), same style as on "2021-Q4" -------------
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$wsQ4.Range("A2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top summary row for "2022-Q1"
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("B2:D2").Style = "Normal"
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.26

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$wsTotal.Range("A2").Value = 0
